$d = $word.ActiveDocument

# --- Fix subject-verb agreement: "needs" -> "need" ------------------------
# A temporary scaffolding bookmark is dropped at the boundary right before
# " ha" (the start of the next run) so the text edit below - which only
# touches the first run - can't ripple forward and coalesce the following,
# untouched runs together. After the edit we delete the scaffolding
# bookmark again; the run split it forced stays in place.
$full = $d.Content.Text
$idxPara = $full.IndexOf("Now consider your own letter")
$idxHa = $full.IndexOf(" ha", $idxPara)
$rBoundary = $d.Range($idxHa, $idxHa)
$d.Bookmarks.Add("_TempSplit", $rBoundary)

$rFix = $d.Content
$rFix.Find.Execute("needs")
$rFix.Text = "need"

$d.Bookmarks("_TempSplit").Delete()

# --- Drop the _GoBack bookmark right after "Write a" -----------------------
# This splits that run the same way Word leaves its last-edit marker.
$rGoBack = $d.Content
$rGoBack.Find.Execute("Write a")
$rGoBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rGoBack)
